$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.876.06"
$ws.Range("D3").Value = "2.649.29"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.03%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "609.43"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "147.84"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +3.02%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("E10").Value = "  +6.84%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "5.57"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.51%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.152"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.62%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "27.53"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "3.121.85"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "63.698.02"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "2.637.38"
$ws.Range("E17").Value = "  -1.13%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "11.75"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +2.75%  "
$ws.Range("E19").Value = "  +4.20%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "346.24"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("E23").Value = "  -3.24%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "66.50"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.57%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "1.65"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +5.48%  "
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("E27").Value = "  +6.54%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "564.53"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +4.52%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "8.21"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +4.22%  "
$ws.Range("E30").Value = "  -1.68%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.996"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.40%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.05"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("E33").Value = "  +5.64%  "
$ws.Range("E34").Value = "  -2.17%  "
$ws.Range("E35").Value = "  +4.30%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "168.61"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -2.28%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  -0.07%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.94"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +6.16%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "19.18"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("E41").Value = "  +0.10%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "165.88"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -4.92%  "
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("E44").Value = "  +1.50%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "22.06"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("E48").Value = "  +16.46%  "
$ws.Range("E49").Value = "  +1.99%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0960"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.45%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "18.86"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.05%  "
